# Update existing biosample data:
#  - harvester (B) renamed from "Retrofitted_2002" to "S.GISH"
#  - experimentDesign (D) filled in with "90minuteInduction"
#  - strain (F) filled in per-genotype (G) with the matching strain name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$strainByGenotype = @{
    "CNAG_00000"             = "KN99alpha"
    "CNAG_00068.CNAG_00156"  = "TDY1960"
    "CNAG_00156"             = "TDY1700"
    "CNAG_06156"             = "TDY1936"
    "CNAG_00068"             = "TDY1945"
    "CNAG_02305"             = "TDY1974"
    "CNAG_02364"             = "TDY2004"
    "CNAG_02305.CNAG_03115"  = "TDY2017"
    "CNAG_03115.CNAG_06252"  = "TDY1963"
    "CNAG_00871"             = "TDY1174"
    "CNAG_06252"             = "TDY1338"
    "CNAG_03115"             = "TDY1948"
    "CNAG_00871.CNAG_06134"  = "TDY1981"
}

for ($r = 2; $r -le 49; $r++) {
    $genotype = $ws.Cells.Item($r, 7).Value()

    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
    $ws.Cells.Item($r, 6).Value = $strainByGenotype[$genotype]
}

$ws.Range("F3:F4").Select()
